# Running only watchlist module
# Set Runmode column (C) so that only the "Watchlist module" (E Suite, row 6)
# is set to run ("Y"); all other modules are turned off ("N").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C2").Value = "N"   # A Suite - IAM module
$ws.Range("C3").Value = "N"   # B Suite - Search module
$ws.Range("C4").Value = "N"   # C Suite - Authoring module
$ws.Range("C5").Value = "N"   # D Suite - Profile module
$ws.Range("C6").Value = "Y"   # E Suite - Watchlist module
$ws.Range("C7").Value = "N"   # F Suite - Notification module

$ws.Range("C6").Select()
